$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value2 = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "65.562.38"
Set-TextValue "E2" "  +3.59%  "

Set-TextValue "D3" "3.491.91"
Set-TextValue "E3" "  +2.36%  "

Set-TextValue "E4" "  +0.02%  "

Set-TextValue "D5" "582.31"
Set-TextValue "E5" "  +2.50%  "

Set-TextValue "D6" "160.94"
Set-TextValue "E6" "  +3.50%  "

Set-TextValue "B7" "XRP"
Set-TextValue "C7" "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D7" "0.613"
Set-TextValue "E7" "  +12.98%  "

Set-TextValue "B8" "USDC"
Set-TextValue "C8" "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D8" "1.00"
Set-TextValue "E8" "  -0.01%  "

Set-TextValue "D9" "3.493.30"
Set-TextValue "E9" "  +2.41%  "

Set-TextValue "D10" "7.31"
Set-TextValue "E10" "  -1.18%  "

Set-TextValue "E11" "  +3.11%  "

Set-TextValue "E12" "  +3.15%  "

Set-TextValue "D13" "4.093.37"
Set-TextValue "E13" "  +2.42%  "

Set-TextValue "E14" "  +0.79%  "

Set-TextValue "E15" "  +3.10%  "

Set-TextValue "D16" "28.66"
Set-TextValue "E16" "  +6.04%  "

Set-TextValue "D17" "65.570.63"
Set-TextValue "E17" "  +3.39%  "

Set-TextValue "D18" "3.494.41"
Set-TextValue "E18" "  +3.74%  "

Set-TextValue "E19" "  +3.36%  "

Set-TextValue "D20" "14.36"
Set-TextValue "E20" "  +1.85%  "

Set-TextValue "D21" "388.95"
Set-TextValue "E21" "  +1.07%  "

Set-TextValue "E22" "  +2.10%  "

Set-TextValue "D23" "0.555"
Set-TextValue "E23" "  +4.16%  "

Set-TextValue "D24" "73.41"
Set-TextValue "E24" "  +2.49%  "

Set-TextValue "D25" "0.997"
Set-TextValue "E25" "  -0.35%  "

Set-TextValue "D26" "0.0000124"
Set-TextValue "E26" "  +5.88%  "

Set-TextValue "D27" "10.10"
Set-TextValue "E27" "  +7.67%  "

Set-TextValue "E28" "  +1.61%  "

Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +0.13%  "

Set-TextValue "B30" "Fetch.AI"
Set-TextValue "C30" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D30" "1.47"
Set-TextValue "E30" "  +9.99%  "

Set-TextValue "B31" "NEARProtocol"
Set-TextValue "C31" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D31" "6.26"
Set-TextValue "E31" "  +4.01%  "

Set-TextValue "E32" "  +3.65%  "

Set-TextValue "D33" "23.71"
Set-TextValue "E33" "  +2.50%  "

Set-TextValue "D34" "7.31"
Set-TextValue "E34" "  +7.44%  "

Set-TextValue "D35" "1.57"
Set-TextValue "E35" "  +8.04%  "

Set-TextValue "D36" "162.73"
Set-TextValue "E36" "  +3.21%  "

Set-TextValue "E37" "  +6.62%  "

Set-TextValue "D38" "3.014.58"
Set-TextValue "E38" "  +4.51%  "

Set-TextValue "E39" "  +1.68%  "

Set-TextValue "D40" "27.41"
Set-TextValue "E40" "  +2.11%  "

Set-TextValue "E41" "  +3.25%  "

Set-TextValue "E42" "  +5.46%  "

Set-TextValue "E43" "  +1.74%  "

Set-TextValue "D44" "42.91"
Set-TextValue "E44" "  +4.45%  "

Set-TextValue "D45" "0.781"
Set-TextValue "E45" "  +3.06%  "

Set-TextValue "D46" "25.31"
Set-TextValue "E46" "  +7.63%  "

Set-TextValue "E47" "  +4.39%  "

Set-TextValue "D48" "321.89"
Set-TextValue "E48" "  +11.41%  "

Set-TextValue "B49" "Stellar"
Set-TextValue "C49" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D49" "0.110"
Set-TextValue "E49" "  +6.93%  "

Set-TextValue "E50" "  +2.58%  "

Set-TextValue "B51" "Cosmos"
Set-TextValue "C51" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D51" "6.75"
Set-TextValue "E51" "  +5.31%  "

